$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.310.58'
$ws.Range("E2").Value = '  +0.57%  '

$ws.Range("D3").Value = '1.667.04'
$ws.Range("E3").Value = '  +1.45%  '

$ws.Range("E4").Value = '  +0.21%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '312.10'
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Range("E5").Value = '  +1.32%  '

$ws.Range("E6").Value = '  +0.17%  '

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.3961'
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Range("E7").Value = '  +1.48%  '

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.3933'
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Range("E8").Value = '  +1.74%  '

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '51.90'
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Range("E9").Value = '  +4.32%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '1.385'
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Range("E10").Value = '  +2.31%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '1.002'
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Range("E11").Value = '  +0.22%  '

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.08569'
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Range("E12").Value = '  -1.28%  '

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '24.38'
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Range("E13").Value = '  +3.01%  '

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '7.311'
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Range("E14").Value = '  +2.69%  '

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '8.008'
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Range("E15").Value = '  +7.37%  '

$ws.Range("E16").Value = '  +3.14%  '

$ws.Range("D17").Value = '1.666.98'
$ws.Range("E17").Value = '  +0.28%  '

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '95.67'
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Range("E18").Value = '  +0.75%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '0.07006'
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Range("E19").Value = '  +1.57%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '20.49'
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Range("E20").Value = '  -0.17%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '6.989'
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Range("E21").Value = '  +1.33%  '

$ws.Range("E22").Value = '  +0.18%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '13.89'
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Range("E23").Value = '  +2.20%  '

$ws.Range("D24").Value = '24.322.31'
$ws.Range("E24").Value = '  +0.65%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '2.531'
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Range("E25").Value = '  +8.57%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '3.104'
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Range("E26").Value = '  +11.46%  '

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '22.48'
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Range("E27").Value = '  +0.43%  '

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '156.85'
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Range("E28").Value = '  -0.52%  '

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '142.38'
$ws.Cells.Item(29, 4).Style = "Normal"
$ws.Range("E29").Value = '  +1.42%  '

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '5.365'
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Range("E30").Value = '  +0.12%  '

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '8.025'
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Range("E31").Value = '  -5.36%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '2.532'
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Range("E32").Value = '  +4.99%  '

$ws.Range("D33").Value = '1.849.70'
$ws.Range("E33").Value = '  +4.62%  '

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '1.060'
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Range("E34").Value = '  +11.65%  '

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.03086'
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Range("E35").Value = '  +6.41%  '

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.08279'
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Range("E36").Value = '  +2.65%  '

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '6.868'
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Range("E37").Value = '  -1.49%  '

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '11.14'
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Range("E38").Value = '  +10.99%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.2757'
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Range("E39").Value = '  +2.83%  '

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.09279'
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Range("E40").Value = '  +0.70%  '

$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.7676'
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Range("E41").Value = '  +1.73%  '

$ws.Range("B42").Value = 'Aptos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '13.75'
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Range("E42").Value = '  +5.59%  '

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '1.442'
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Range("E43").Value = '  -1.06%  '

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '16.64'
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Range("E44").Value = '  +4.09%  '

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.7066'
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Range("E45").Value = '  +2.33%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '2.532'
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Range("E46").Value = '  +2.84%  '

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '4.125'
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Range("E47").Value = '  +0.88%  '

$ws.Range("E48").Value = '  +0.25%  '

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.08415'
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Range("E49").Value = '  +0.07%  '

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '136.34'
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Range("E50").Value = '  +1.86%  '

$ws.Range("E51").Value = '  -0.05%  '
